# Update market price / profit data cells across multiple worksheets
# as refreshed by the scheduled Universalis-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H33").Value = 164.17857
$ws.Range("I33").Value = 70.59090999999999
$ws.Range("J33").Value = 507.33334
$ws.Range("K33").Value = 70.59090999999999
$ws.Range("L33").Value = 507.33334
$ws.Range("M33").Value = 158.40909
$ws.Range("N33").Value = -965.33334

$ws.Range("H135").Value = 1349.0869
$ws.Range("I135").Value = 1707.2
$ws.Range("J135").Value = 677.625
$ws.Range("K135").Value = 15364.8
$ws.Range("L135").Value = 6098.625
$ws.Range("M135").Value = -12829.8
$ws.Range("N135").Value = -11168.625

$ws = $wb.Worksheets("ARM")
$ws.Range("H61").Value = 1413.8372
$ws.Range("I61").Value = 1453.2894
$ws.Range("J61").Value = 1114
$ws.Range("K61").Value = 1453.2894
$ws.Range("L61").Value = 1114
$ws.Range("M61").Value = -1241.2894
$ws.Range("N61").Value = -1538

$ws.Range("H74").Value = 1224.2572
$ws.Range("I74").Value = 1191.4517
$ws.Range("J74").Value = 1478.5
$ws.Range("K74").Value = 1191.4517
$ws.Range("L74").Value = 1478.5
$ws.Range("M74").Value = -317.4517000000001
$ws.Range("N74").Value = -3226.5

$ws.Range("H77").Value = 1224.2572
$ws.Range("I77").Value = 1191.4517
$ws.Range("J77").Value = 1478.5
$ws.Range("K77").Value = 5957.2585
$ws.Range("L77").Value = 7392.5
$ws.Range("M77").Value = -1589.2585
$ws.Range("N77").Value = -16128.5

$ws.Range("H97").Value = 4182.85
$ws.Range("I97").Value = 5236.4287
$ws.Range("J97").Value = 1724.5
$ws.Range("K97").Value = 5236.4287
$ws.Range("L97").Value = 1724.5
$ws.Range("M97").Value = -4740.4287
$ws.Range("N97").Value = -2716.5

$ws.Range("H132").Value = 1248.92
$ws.Range("I132").Value = 1073.3617
$ws.Range("K132").Value = 3220.0851
$ws.Range("M132").Value = -690.0850999999998

$ws.Range("H136").Value = 1413.8372
$ws.Range("I136").Value = 1453.2894
$ws.Range("J136").Value = 1114
$ws.Range("K136").Value = 4359.8682
$ws.Range("L136").Value = 3342
$ws.Range("M136").Value = -1809.8682
$ws.Range("N136").Value = -8442

$ws = $wb.Worksheets("BSM")
$ws.Range("H134").Value = 30830.943
$ws.Range("I134").Value = 2041.3871
$ws.Range("J134").Value = 253950
$ws.Range("K134").Value = 6124.1613
$ws.Range("L134").Value = 761850
$ws.Range("M134").Value = -3589.1613
$ws.Range("N134").Value = -766920

$ws = $wb.Worksheets("CRP")
$ws.Range("H31").Value = 43472.44
$ws.Range("I31").Value = 3825.6875
$ws.Range("J31").Value = 113955.555
$ws.Range("K31").Value = 3825.6875
$ws.Range("L31").Value = 113955.555
$ws.Range("M31").Value = -3530.6875
$ws.Range("N31").Value = -114545.555

$ws.Range("H34").Value = 43472.44
$ws.Range("I34").Value = 3825.6875
$ws.Range("J34").Value = 113955.555
$ws.Range("K34").Value = 3825.6875
$ws.Range("L34").Value = 113955.555
$ws.Range("M34").Value = -3623.6875
$ws.Range("N34").Value = -114359.555

$ws = $wb.Worksheets("CUL")
$ws.Range("H29").Value = 188.33333
$ws.Range("J29").Value = 232.5
$ws.Range("L29").Value = 697.5
$ws.Range("N29").Value = -1251.5

$ws.Range("H40").Value = 135
$ws.Range("I40").Value = 176.66667
$ws.Range("J40").Value = 93.333336
$ws.Range("K40").Value = 706.66668
$ws.Range("L40").Value = 373.333344
$ws.Range("M40").Value = -637.66668
$ws.Range("N40").Value = -511.333344

$ws.Range("H87").Value = 13799.8
$ws.Range("I87").Value = 5399.6
$ws.Range("J87").Value = 17999.9
$ws.Range("K87").Value = 16198.8
$ws.Range("L87").Value = 53999.7
$ws.Range("M87").Value = -14950.8
$ws.Range("N87").Value = -56495.7

$ws.Range("H90").Value = 13799.8
$ws.Range("I90").Value = 5399.6
$ws.Range("J90").Value = 17999.9
$ws.Range("K90").Value = 48596.4
$ws.Range("L90").Value = 161999.1
$ws.Range("M90").Value = -42356.4
$ws.Range("N90").Value = -174479.1

$ws.Range("H99").Value = 144999.86
$ws.Range("I99").Value = 250749.75
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 752249.25
$ws.Range("L99").Value = 12000
$ws.Range("M99").Value = -750003.25
$ws.Range("N99").Value = -16492

$ws.Range("H131").Value = 17277212
$ws.Range("I131").Value = 167000160
$ws.Range("J131").Value = 1486.7307
$ws.Range("K131").Value = 501000480
$ws.Range("L131").Value = 4460.1921
$ws.Range("M131").Value = -500995440
$ws.Range("N131").Value = -14540.1921

$ws.Range("H139").Value = 39588.77
$ws.Range("I139").Value = 40932.32
$ws.Range("J139").Value = 6000
$ws.Range("K139").Value = 122796.96
$ws.Range("L139").Value = 18000
$ws.Range("M139").Value = -117656.96
$ws.Range("N139").Value = -28280

$ws = $wb.Worksheets("GSM")
$ws.Range("H141").Value = 56366.125
$ws.Range("J141").Value = 56366.125
$ws.Range("L141").Value = 56366.125
$ws.Range("N141").Value = -66726.125

$ws = $wb.Worksheets("LTW")
$ws.Range("H55").Value = 303.2143
$ws.Range("I55").Value = 433.16666
$ws.Range("J55").Value = 205.75
$ws.Range("K55").Value = 433.16666
$ws.Range("L55").Value = 205.75
$ws.Range("M55").Value = -260.16666
$ws.Range("N55").Value = -551.75

$ws.Range("H132").Value = 2181.889
$ws.Range("I132").Value = 2152.9
$ws.Range("J132").Value = 2239.8667
$ws.Range("K132").Value = 6458.700000000001
$ws.Range("L132").Value = 6719.6001
$ws.Range("M132").Value = -3928.700000000001
$ws.Range("N132").Value = -11779.6001

$ws.Range("H136").Value = 3859.5217
$ws.Range("I136").Value = 1992.2941
$ws.Range("J136").Value = 9150
$ws.Range("K136").Value = 5976.8823
$ws.Range("L136").Value = 27450
$ws.Range("M136").Value = -3426.8823
$ws.Range("N136").Value = -32550

$ws = $wb.Worksheets("WVR")
$ws.Range("H132").Value = 813.3469
$ws.Range("I132").Value = 691.0227
$ws.Range("J132").Value = 1889.8
$ws.Range("K132").Value = 2073.0681
$ws.Range("L132").Value = 5669.4
$ws.Range("M132").Value = 456.9319
$ws.Range("N132").Value = -10729.4

$ws.Range("H136").Value = 3057.0833
$ws.Range("I136").Value = 3280
$ws.Range("J136").Value = 1942.5
$ws.Range("K136").Value = 9840
$ws.Range("L136").Value = 5827.5
$ws.Range("M136").Value = -7290
$ws.Range("N136").Value = -10927.5
